# Update "Förändrad" (Changed) date column (C) for all data rows.
# The diff shows every data row's column C value changing from 45179 to 45180
# (i.e. the date serial advances by one day, 2023-09-10 -> 2023-09-11),
# for rows 2 through 132.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 132 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
